$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (CANALINI / LUCA / MAIA): fill in skill ratings that were blank ---
$ws.Range("D2").Value = "1 - Beginner"
$ws.Range("E2").Value = "3 - Good"
$ws.Range("F2").Value = "1 - Beginner"
$ws.Range("G2").Value = "1 - Beginner"
$ws.Range("H2").Value = "1 - Beginner"
$ws.Range("I2").Value = "3 - Good"

# --- Row 11 (GHIMIRE / PAMIR / MSCV): fill in skill ratings + remark ---
# (written before row 10's remark so the shared-string table picks up the
# same ordering/indices as the target workbook)
$ws.Range("D11").Value = "2 - Average"
$ws.Range("E11").Value = "4 - Excellent"
$ws.Range("F11").Value = "1 - Beginner"
$ws.Range("G11").Value = "3 - Good"
$ws.Range("H11").Value = "2 - Average"
$ws.Range("I11").Value = "3 - Good"
$ws.Range("J11").Value = "Have experience with research, could do sketches in documentation"

# --- Row 10 (DOUSAI / NAYEE MUDDIN KHAN / MSCV): update ratings + add remark ---
$ws.Range("D10").Value = "4 - Excellent"
$ws.Range("F10").Value = "2 - Average"
$ws.Range("H10").Value = "2 - Average"
$ws.Range("J10").Value = "Worked in real software building environment and managed a startup of about 12 collegues as co founder"

# --- Row heights (wrapped remark text grew taller on re-save) ---
$ws.Rows.Item(3).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 180
$ws.Rows.Item(8).RowHeight = 90
$ws.Rows.Item(10).RowHeight = 195
$ws.Rows.Item(11).RowHeight = 135
$ws.Rows.Item(12).RowHeight = 75

# --- Selection moved to I2 on the resave ---
[void]$ws.Range("I2").Select()
